# Append the latest Adafruit IO reading as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 73

# Column C holds a numeric-looking reading ("25") but the sheet stores every
# value as text, so force a text number format before assigning it -
# otherwise Excel would coerce it into a real number.
$ws.Cells.Item($row, 3).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
